$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add new header cells P1=14, Q1=15, using the same style as the
#     rest of the header row (style index 1: bold, centered, bordered). ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Rows 2-25: swap values in columns I/K and M/O (1<->2), and add new
#     columns P and Q, each filled with 2 (no special style, matching the
#     rest of the data columns). ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q = 2
}
